$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Column D: "Avaliação" (stars given to each book) ---
$ws.Range("D1").Value = "Avaliação"
$ws.Range("D1").Style = $ws.Range("A1").Style

$avaliacoes = @(2,5,3,4,4,1,1,4,5,1,2,4,1,2,3,1,5,1,5,3)
for ($i = 0; $i -lt $avaliacoes.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 4).Value = $avaliacoes[$i]
}

# --- KPI block ---
$ws.Range("A23").Value = "Indicadores de Performance"

$ws.Range("A24").Value = "Percentual Bem Avaliados (%)"
$ws.Range("B24").Value = 40

$ws.Range("A25").Value = "Percentual Estoque Crítico (%)"
$ws.Range("B25").Value = 0

$ws.Range("A26").Value = "Preço Médio Bem Avaliados (£)"
$ws.Range("B26").Value = 35.43

# --- Ratings distribution table ---
$ws.Range("A27").Value = "Avaliação"
$ws.Range("B27").Value = "Contagem"
$ws.Range("C27").Value = "Percentual"

$ratings = @(
    @("1 estrelas", 6, 30),
    @("2 estrelas", 3, 15),
    @("3 estrelas", 3, 15),
    @("4 estrelas", 4, 20),
    @("5 estrelas", 4, 20)
)

for ($i = 0; $i -lt $ratings.Count; $i++) {
    $row = 28 + $i
    $ws.Cells.Item($row, 1).Value = $ratings[$i][0]
    $ws.Cells.Item($row, 2).Value = $ratings[$i][1]
    $ws.Cells.Item($row, 3).Value = $ratings[$i][2]
}
